# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
#
# This "Estado de Cuenta" sheet lists the periods in arrears for a single
# worker. The worker changes (new document number / name), the list of
# overdue periods grows from 4 to 5 rows (now listed newest-period-first),
# and the summary totals (VALOR MORA / Cant. Periodos) are refreshed to
# match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Grow the worker's detail table from 4 to 5 rows -----------------------
# Row 19 currently carries the special "closing" border style (it's the
# last row of the table). Insert a fresh row above it so that style stays
# at the bottom of the (now 5-row) table.
$ws.Rows(20).Insert()

# The inserted row20 starts out blank with a blended/default look; pull the
# exact "closing row" formatting down from the row it displaced (still row19
# at this point).
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)   # xlPasteFormats

# Row 19 itself becomes a normal interior row now, so give it the same
# formatting as the row above it (row18, a normal row).
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Update the account summary block --------------------------------------
$ws.Range("E11").Value = 284700   # VALOR MORA total
$ws.Range("F13").Value = 5        # Cant. Periodos

# --- Replace the worker identity on every detail row ------------------------
$rows = 16, 17, 18, 19, 20
foreach ($r in $rows) {
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "1002059412"
    $ws.Range("D$r").Value = "ERIKA PATRICIA PEREZ MARTINEZ"
    $ws.Range("F$r").Value = 56940
    $ws.Range("G$r").Value = 1423500
}

# --- Refresh the overdue periods, newest first ------------------------------
$ws.Range("E16").Value = "2507"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2505"
$ws.Range("E19").Value = "2504"
$ws.Range("E20").Value = "2503"
